$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# ValidNames (sheet2): append two more valid names
# ---------------------------------------------------------------------
$wsValidNames = $wb.Worksheets.Item("ValidNames")
$wsValidNames.Range("A5").Value = "Srihitha"
$wsValidNames.Range("A6").Value = "Rashmi"
$wsValidNames.Columns.Item(1).ColumnWidth = 20.6
[void]$wsValidNames.Range("A7").Select()

# ---------------------------------------------------------------------
# InvalidNames (sheet3): append four more invalid-name samples
# ---------------------------------------------------------------------
$wsInvalidNames = $wb.Worksheets.Item("InvalidNames")
$wsInvalidNames.Range("A5").Value = 123
$wsInvalidNames.Range("A6").Value = "…"
# Leading apostrophe forces the text quote-prefix style (matches "--" in A3)
$wsInvalidNames.Range("A7").Value = "'++++++"
$wsInvalidNames.Range("A8").Value = 98400021
$wsInvalidNames.Columns.Item(1).ColumnWidth = 8.25
[void]$wsInvalidNames.Range("A9").Select()

# ---------------------------------------------------------------------
# ValidEmails (sheet4): append five more valid emails, each a mailto hyperlink
# ---------------------------------------------------------------------
$wsValidEmails = $wb.Worksheets.Item("ValidEmails")

$wsValidEmails.Range("A4").Value = "a7@gmail.com"
[void]$wsValidEmails.Hyperlinks.Add($wsValidEmails.Range("A4"), "mailto:a7@gmail.com")

$wsValidEmails.Range("A5").Value = "a7@yahoo.com"
[void]$wsValidEmails.Hyperlinks.Add($wsValidEmails.Range("A5"), "mailto:a7@yahoo.com")

$wsValidEmails.Range("A6").Value = "cgi@cgi.com"
[void]$wsValidEmails.Hyperlinks.Add($wsValidEmails.Range("A6"), "mailto:cgi@cgi.com")

$wsValidEmails.Range("A7").Value = "a7@cgi.com"
[void]$wsValidEmails.Hyperlinks.Add($wsValidEmails.Range("A7"), "mailto:a7@cgi.com")

$wsValidEmails.Range("A8").Value = "A7TEAM@GAMAIL.COM"
[void]$wsValidEmails.Hyperlinks.Add($wsValidEmails.Range("A8"), "mailto:A7TEAM@GAMAIL.COM")

# Hyperlinks.Add creates a fresh cellXf for the hyperlink font; re-point these
# cells back at the workbook's existing "Hyperlink" cell style so the new
# cells land on the same style index as the pre-existing hyperlink (A2).
$wsValidEmails.Range("A4:A8").Style = "Hyperlink"

[void]$wsValidEmails.Range("A9").Select()

# ---------------------------------------------------------------------
# InvalidEmails (sheet5): append four more invalid-email samples, and this
# becomes the active sheet/tab
# ---------------------------------------------------------------------
$wsInvalidEmails = $wb.Worksheets.Item("InvalidEmails")
$wsInvalidEmails.Range("A5").Value = "CGIDOTCOM"
$wsInvalidEmails.Range("A6").Value = "cgidotcom"
$wsInvalidEmails.Range("A7").Value = "a7atcgidotcom"
$wsInvalidEmails.Range("A8").Value = 870022
$wsInvalidEmails.Columns.Item(1).ColumnWidth = 26.1
[void]$wsInvalidEmails.Range("A12").Select()
[void]$wsInvalidEmails.Activate()

# ---------------------------------------------------------------------
# Sheet1: no data changes, just widen column C and it's no longer the
# selected tab (InvalidEmails is, per above)
# ---------------------------------------------------------------------
$wsSheet1 = $wb.Worksheets.Item("Sheet1")
$wsSheet1.Columns.Item(3).ColumnWidth = 15.75
